# Update stack-trace line numbers and trailing frames in the error-message run
# (Moving from 3.1.1 to 3.2.0).
$d = $word.ActiveDocument

$d.Content.Find.Execute("M2DocEvaluator.java:1480", $true, $false, $false, $false, $false, $true, 1, $false, "M2DocEvaluator.java:1569", 2) | Out-Null
$d.Content.Find.Execute("M2DocEvaluator.java:1242", $true, $false, $false, $false, $false, $true, 1, $false, "M2DocEvaluator.java:1331", 2) | Out-Null
$d.Content.Find.Execute("M2DocEvaluator.java:1467", $true, $false, $false, $false, $false, $true, 1, $false, "M2DocEvaluator.java:1556", 2) | Out-Null
$d.Content.Find.Execute("M2DocEvaluator.java:1491", $true, $false, $false, $false, $false, $true, 1, $false, "M2DocEvaluator.java:1580", 2) | Out-Null
$d.Content.Find.Execute("M2DocEvaluator.java:297)", $true, $false, $false, $false, $false, $true, 1, $false, "M2DocEvaluator.java:301)", 2) | Out-Null
$d.Content.Find.Execute("M2DocEvaluator.java:282)", $true, $false, $false, $false, $false, $true, 1, $false, "M2DocEvaluator.java:286)", 2) | Out-Null
$d.Content.Find.Execute("M2DocUtils.java:845", $true, $false, $false, $false, $false, $true, 1, $false, "M2DocUtils.java:853", 2) | Out-Null
$d.Content.Find.Execute("AbstractTemplatesTestSuite.java:514", $true, $false, $false, $false, $false, $true, 1, $false, "AbstractTemplatesTestSuite.java:518", 2) | Out-Null
$d.Content.Find.Execute("AbstractTemplatesTestSuite.java:421", $true, $false, $false, $false, $false, $true, 1, $false, "AbstractTemplatesTestSuite.java:414", 2) | Out-Null
$d.Content.Find.Execute("GeneratedMethodAccessor73", $true, $false, $false, $false, $false, $true, 1, $false, "GeneratedMethodAccessor5", 2) | Out-Null

# Replace the trailing JDT/Eclipse test-runner frames with the Maven Surefire / Tycho / Equinox frames
$oldBlock = "`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)"
$newBlock = "`tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:365)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeWithRerun(JUnit4Provider.java:273)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:238)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:159)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:206)`n`tat org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:161)`n`tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:84)`n`tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:113)`n`tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:593)`n`tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:205)`n`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:137)`n`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:107)`n`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:401)`n`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:255)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:657)`n`tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:594)`n`tat org.eclipse.equinox.launcher.Main.run(Main.java:1447)`n`tat org.eclipse.equinox.launcher.Main.main(Main.java:1420)"
$d.Content.Find.Execute($oldBlock, $true, $false, $false, $false, $false, $true, 1, $false, $newBlock, 2) | Out-Null

$d.Save()
